$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 127: clear the "NA" placeholder text that used to live in C127
$ws.Range("C127").ClearContents()

# New row 128
$ws.Range("A128").Value = "'2025-06-11"
$ws.Range("B128").Value = "ruissellement"
$ws.Range("C128").Value = 113
$ws.Range("D128").Value = 2

# New row 129
$ws.Range("A129").Value = "'2025-06-11"
$ws.Range("B129").Value = "ruissellement"
$ws.Range("C129").Value = 115
$ws.Range("D129").Value = 2

# New row 130
$ws.Range("A130").Value = "'2025-06-11"
$ws.Range("B130").Value = "bonnes pratiques"
$ws.Range("C130").Value = 116
$ws.Range("D130").Value = 1
